# Update ATT table figures (Std_Error / CI_Lower_95 / CI_Upper_95 columns)
# for "Overall ATT (Group aggregation)", "Overall ATT (Dynamic aggregation)",
# "Cohort 2014" and "Cohort 2015" rows.

$d = $word.ActiveDocument

# Row: Overall ATT (Group aggregation)
$d.Content.Find.Execute("1.80", $true, $false, $false, $false, $false, $true, 1, $false, "1.83", 2)
$d.Content.Find.Execute("-17.57", $true, $false, $false, $false, $false, $true, 1, $false, "-17.62", 2)
$d.Content.Find.Execute("-10.51", $true, $false, $false, $false, $false, $true, 1, $false, "-10.46", 2)

# Row: Overall ATT (Dynamic aggregation)
$d.Content.Find.Execute("1.98", $true, $false, $false, $false, $false, $true, 1, $false, "1.99", 2)
$d.Content.Find.Execute("-17.92", $true, $false, $false, $false, $false, $true, 1, $false, "-17.95", 2)
$d.Content.Find.Execute("-10.16", $true, $false, $false, $false, $false, $true, 1, $false, "-10.13", 2)

# Row: Cohort 2014
$d.Content.Find.Execute("2.51", $true, $false, $false, $false, $false, $true, 1, $false, "2.47", 2)
$d.Content.Find.Execute("-21.73", $true, $false, $false, $false, $false, $true, 1, $false, "-21.65", 2)
$d.Content.Find.Execute("-11.89", $true, $false, $false, $false, $false, $true, 1, $false, "-11.97", 2)

# Row: Cohort 2015
$d.Content.Find.Execute("2.82", $true, $false, $false, $false, $false, $true, 1, $false, "2.92", 2)
$d.Content.Find.Execute("-17.25", $true, $false, $false, $false, $false, $true, 1, $false, "-17.45", 2)
$d.Content.Find.Execute("-6.21", $true, $false, $false, $false, $false, $true, 1, $false, "-6.01", 2)
